# revised based on Meng's comments
# The "funding" sheet's task-rank column (F) used "负责人1" (with assorted
# whitespace variants) for the first six projects; drop the trailing "1" so
# it simply reads "负责人".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("funding")

$ws.Range("F2").Value = "负责人"
$ws.Range("F3").Value = "负责人"
$ws.Range("F4").Value = "负责人"
$ws.Range("F5").Value = "负责人"
$ws.Range("F6").Value = "负责人"
$ws.Range("F7").Value = "负责人"

$ws.Range("H7").Select()
